# Apply the gene-id text shortening + selection change described in the
# commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shorten the gene identifiers in column A (rows 2-4).
$ws.Range("A2").Value = "G00001"
$ws.Range("A3").Value = "G00002"
$ws.Range("A4").Value = "G00003"

# Move the saved selection from D5 to A5.
$ws.Range("A5").Select()
